$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price/Volume columns so that numeric-looking
# strings (e.g. "530.24", "17.50") are written verbatim as text instead of
# being auto-converted to numbers by Excel's smart-entry parser.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.613.70"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "3.345.19"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "530.24"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").Value = "172.18"
$ws.Range("E6").Value = "  -5.00%  "
$ws.Range("D7").Value = "0.594"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "3.345.28"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").Value = "52.94"
$ws.Range("E11").Value = "  -7.78%  "
$ws.Range("E12").Value = "  +2.45%  "
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "3.884.36"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "3.345.68"
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "17.50"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "63.482.81"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("D22").Value = "372.08"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "11.35"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "81.56"
$ws.Range("E24").Value = "  +2.05%  "
$ws.Range("D25").Value = "3.72"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +5.43%  "
$ws.Range("D27").Value = "6.19"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").Value = "2.70"
$ws.Range("E28").Value = "  +3.02%  "
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "8.23"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "635.76"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").Value = "11.18"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  -2.07%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "36.28"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").Value = "0.0₃0721"
$ws.Range("E40").Value = "  +10.71%  "
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  +6.34%  "
$ws.Range("D43").Value = "2.947.21"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +7.55%  "
$ws.Range("D46").Value = "2.70"
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").Value = "0.0396"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").Value = "3.09"
$ws.Range("E48").Value = "  +5.86%  "
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "135.95"
$ws.Range("E51").Value = "  +4.38%  "

# Restore the original (General) number format / default style so the
# cells match their pre-edit appearance.
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("E2:E51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
